# Update 2025 (row 9) sales figures on Sheet1 as part of
# "terminei vendas BIBI e arrumei vendas ADD na analise e no dash"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 2266038.68
$ws.Range("C9").Value = 325971.99
$ws.Range("D9").Value = 2592010.67
$ws.Range("E9").Value = 12.57602809173621
$ws.Range("F9").Value = 87.42397190826379
$ws.Range("G9").Value = -68.49643253777445
$ws.Range("H9").Value = -59.07847364763155
$ws.Range("I9").Value = -60.56121065787934
$ws.Range("J9").Value = 22166
$ws.Range("K9").Value = 938
$ws.Range("L9").Value = 23104
